$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update header row (A1:D1) with new snake_case column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case fix for municipality/state names containing lowercase connector words
#    (de, del, la, las, los, el, y) -> capitalized, when not the first word
$stringFixes = @(
    @('B8', 'Pabellón De Arteaga'),
    @('B9', 'Rincón De Romos'),
    @('B10', 'San Francisco De Los Romo'),
    @('B15', 'Playas De Rosarito'),
    @('B31', 'Amatenango De La Frontera'),
    @('B38', 'Chiapa De Corzo'),
    @('B41', 'Comitán De Domínguez'),
    @('B63', 'Salto De Agua'),
    @('B64', 'San Cristóbal De Las Casas'),
    @('B97', 'Coyame Del Sotol'),
    @('B108', 'Guadalupe Y Calvo'),
    @('B111', 'Hidalgo Del Parral'),
    @('B136', 'San Francisco De Borja'),
    @('B137', 'San Francisco De Conchos'),
    @('B138', 'San Francisco Del Oro'),
    @('B146', 'Valle De Zaragoza'),
    @('B173', 'Villa De Álvarez'),
    @('A175', 'Ciudad De México'),
    @('B179', 'Cuajimalpa De Morelos'),
    @('B194', 'Coneto De Comonfort'),
    @('B208', 'Nombre De Dios'),
    @('B212', 'Pánuco De Coronado'),
    @('B219', 'San Juan De Guadalupe'),
    @('B220', 'San Juan Del Río'),
    @('B221', 'San Luis Del Cordero'),
    @('B222', 'San Pedro Del Gallo'),
    @('A232', 'Estado De México'),
    @('B232', 'Acambay De Ruíz Castañeda'),
    @('B235', 'Almoloya De Alquisiras'),
    @('B236', 'Almoloya De Juárez'),
    @('B241', 'Atizapán De Zaragoza'),
    @('B249', 'Coacalco De Berriozábal'),
    @('B255', 'Ecatepec De Morelos'),
    @('B260', 'Ixtapan De La Sal'),
    @('B267', 'Naucalpan De Juárez'),
    @('B273', 'San Felipe Del Progreso'),
    @('B274', 'San Martín De Las Pirámides'),
    @('B276', 'San Simón De Guerrero'),
    @('B283', 'Tenango Del Valle'),
    @('B290', 'Tlalnepantla De Baz'),
    @('B296', 'Valle De Bravo'),
    @('B297', 'Valle De Chalco Solidaridad'),
    @('B298', 'Villa De Allende'),
    @('B306', 'San Miguel De Allende'),
    @('B307', 'Apaseo El Alto'),
    @('B308', 'Apaseo El Grande'),
    @('B316', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B327', 'Purísima Del Rincón'),
    @('B331', 'San Diego De La Unión'),
    @('B333', 'San Francisco Del Rincón'),
    @('B335', 'San Luis De La Paz'),
    @('B336', 'Santa Cruz De Juventino Rosas'),
    @('B337', 'Silao De La Victoria'),
    @('B341', 'Valle De Santiago'),
    @('B347', 'Acapulco De Juárez'),
    @('B348', 'Ajuchitlán Del Progreso'),
    @('B349', 'Alcozauca De Guerrero'),
    @('B353', 'Atlamajalcingo Del Monte'),
    @('B355', 'Atoyac De Álvarez'),
    @('B356', 'Ayutla De Los Libres'),
    @('B359', 'Buenavista De Cuéllar'),
    @('B360', 'Chilapa De Álvarez'),
    @('B361', 'Chilpancingo De Los Bravo'),
    @('B362', 'Coahuayutla De José María Izazaga'),
    @('B366', 'Coyuca De Benítez'),
    @('B367', 'Coyuca De Catalán'),
    @('B371', 'Cuetzala Del Progreso'),
    @('B372', 'Cutzamala De Pinzón'),
    @('B378', 'Huitzuco De Los Figueroa'),
    @('B379', 'Iguala De La Independencia'),
    @('B381', 'Ixcateopan De Cuauhtémoc'),
    @('B382', 'Zihuatanejo De Azueta'),
    @('B384', 'La Unión De Isidoro Montes De Oca'),
    @('B387', 'Mártir De Cuilapan'),
    @('B399', 'Taxco De Alarcón'),
    @('B401', 'Técpan De Galeana'),
    @('B403', 'Tepecoacuilco De Trujano'),
    @('B405', 'Tixtla De Guerrero'),
    @('B408', 'Tlalixtaquilla De Maldonado'),
    @('B409', 'Tlapa De Comonfort'),
    @('B420', 'Atotonilco De Tula'),
    @('B421', 'Atotonilco El Grande'),
    @('B426', 'Cuautepec De Hinojosa'),
    @('B433', 'Jacala De Ledezma'),
    @('B437', 'Mineral Del Chico'),
    @('B438', 'Mixquiahuala De Juárez'),
    @('B439', 'Pachuca De Soto'),
    @('B443', 'Santiago De Anaya'),
    @('B444', 'Santiago Tulantepec De Lugo Guerrero'),
    @('B449', 'Tezontepec De Aldama'),
    @('B453', 'Tula De Allende'),
    @('B454', 'Tulancingo De Bravo'),
    @('B455', 'Zacualtipán De Ángeles'),
    @('B460', 'Acatlán De Juárez'),
    @('B461', 'Ahualulco De Mercado'),
    @('B466', 'Atemajac De Brizuela'),
    @('B469', 'Atotonilco El Alto'),
    @('B470', 'Autlán De Navarro'),
    @('B475', 'Cañadas De Obregón'),
    @('B482', 'Cuautitlán De García Barragán'),
    @('B490', 'Encarnación De Díaz'),
    @('B495', 'Huejuquilla El Alto'),
    @('B496', 'Ixtlahuacán Del Río'),
    @('B499', 'Jilotlán De Los Dolores'),
    @('B504', 'Lagos De Moreno'),
    @('B511', 'Ojuelos De Jalisco'),
    @('B516', 'San Diego De Alejandría'),
    @('B517', 'San Juan De Los Lagos'),
    @('B518', 'San Juanito De Escobedo'),
    @('B520', 'San Martín De Bolaños'),
    @('B522', 'San Miguel El Alto'),
    @('B523', 'San Sebastián Del Oeste'),
    @('B524', 'Santa María De Los Ángeles'),
    @('B527', 'Tamazula De Gordiano'),
    @('B530', 'Teocuitatlán De Corona'),
    @('B531', 'Tepatitlán De Morelos'),
    @('B533', 'Tizapán El Alto'),
    @('B534', 'Tlajomulco De Zúñiga'),
    @('B545', 'Unión De Tula'),
    @('B546', 'Valle De Guadalupe'),
    @('B551', 'Yahualica De González Gallo'),
    @('B552', 'Zacoalco De Torres'),
    @('B555', 'Zapotlán El Grande'),
    @('B574', 'Coalcomán De Vázquez Pallares'),
    @('B576', 'Cojumatlán De Régules'),
    @('B633', 'Tiquicheo De Nicolás Romero'),
    @('B663', 'Puente De Ixtla'),
    @('B666', 'Tetela Del Volcán'),
    @('B667', 'Tlaltizapán De Zapata'),
    @('B677', 'Amatlán De Cañas'),
    @('B678', 'Bahía De Banderas'),
    @('B682', 'Ixtlán Del Río'),
    @('B689', 'Santa María Del Oro'),
    @('B706', 'Mier Y Noriega'),
    @('B710', 'San Nicolás De Los Garza'),
    @('B713', 'Acatlán De Pérez Figueroa'),
    @('B718', 'Coicoyán De Las Flores'),
    @('B721', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B722', 'Heroica Ciudad De Huajuapan De León'),
    @('B723', 'Heroica Ciudad De Tlaxiaco'),
    @('B724', 'Ixtlán De Juárez'),
    @('B725', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B728', 'Mariscala De Juárez'),
    @('B730', 'Miahuatlán De Porfirio Díaz'),
    @('B732', 'Nejapa De Madero'),
    @('B733', 'Oaxaca De Juárez'),
    @('B734', 'Ocotlán De Morelos'),
    @('B735', 'Putla Villa De Guerrero'),
    @('B736', 'Rojas De Cuauhtémoc'),
    @('B752', 'San José Del Peñasco'),
    @('B768', 'San Miguel Del Puerto'),
    @('B783', 'San Pedro Y San Pablo Ayutla'),
    @('B792', 'Santa María Del Tule'),
    @('B803', 'Santiago Del Río'),
    @('B815', 'Santo Domingo De Morelos'),
    @('B821', 'Tamazulápam Del Espíritu Santo'),
    @('B822', 'Tezoatlán De Segura Y Luna'),
    @('B823', 'Tlacolula De Matamoros'),
    @('B824', 'Tlalixtac De Cabrera'),
    @('B825', 'Villa De Tututepec De Melchor Ocampo'),
    @('B826', 'Villa De Zaachila'),
    @('B827', 'Villa Sola De Vega'),
    @('B828', 'Zimatlán De Álvarez'),
    @('B843', 'Chalchicomula De Sesma'),
    @('B852', 'Cuayuca De Andrade'),
    @('B859', 'Huehuetlán El Chico'),
    @('B868', 'Los Reyes De Juárez'),
    @('B874', 'Palmar De Bravo'),
    @('B884', 'San Salvador El Seco'),
    @('B885', 'San Salvador El Verde'),
    @('B887', 'Tecali De Herrera'),
    @('B893', 'Tepanco De López'),
    @('B895', 'Tepexi De Rodríguez'),
    @('B897', 'Tetela De Ocampo'),
    @('B902', 'Tlacotepec De Benito Juárez'),
    @('B909', 'Totoltepec De Guerrero'),
    @('B919', 'Amealco De Bonfil'),
    @('B921', 'Cadereyta De Montes'),
    @('B925', 'Jalpan De Serra'),
    @('B926', 'Landa De Matamoros'),
    @('B929', 'Pinal De Amoles'),
    @('B932', 'San Juan Del Río'),
    @('B941', 'Axtla De Terrazas'),
    @('B945', 'Ciudad Del Maíz'),
    @('B952', 'Mexquitic De Carmona'),
    @('B960', 'Santa María Del Río'),
    @('B967', 'Villa De Arriaga'),
    @('B968', 'Villa De Guadalupe'),
    @('B969', 'Villa De Ramos'),
    @('B970', 'Villa De Reyes'),
    @('B1011', 'Nacozari De García'),
    @('B1028', 'Jalpa De Méndez'),
    @('B1049', 'Soto La Marina'),
    @('B1055', 'Acuamanala De Miguel Hidalgo'),
    @('B1061', 'Muñoz De Domingo Arenas'),
    @('B1062', 'Nanacamilpa De Mariano Arista'),
    @('B1064', 'San Pablo Del Monte'),
    @('B1079', 'Alto Lucero De Gutiérrez Barrios'),
    @('B1082', 'Amatlán De Los Reyes'),
    @('B1090', 'Boca Del Río'),
    @('B1093', 'Castillo De Teayo'),
    @('B1095', 'Cazones De Herrera'),
    @('B1103', 'Cosamaloapan De Carpio'),
    @('B1104', 'Cosautlán De Carvajal'),
    @('B1113', 'Hueyapan De Ocampo'),
    @('B1114', 'Ignacio De La Llave'),
    @('B1117', 'Ixhuatlán Del Café'),
    @('B1118', 'Ixhuatlán Del Sureste'),
    @('B1129', 'Lerdo De Tejada'),
    @('B1132', 'Martínez De La Torre'),
    @('B1134', 'Medellín De Bravo'),
    @('B1137', 'Nanchital De Lázaro Cárdenas Del Río'),
    @('B1145', 'Paso De Ovejas'),
    @('B1146', 'Paso Del Macho'),
    @('B1149', 'Poza Rica De Hidalgo'),
    @('B1156', 'Soledad De Doblado'),
    @('B1173', 'Vega De Alatorre'),
    @('B1197', 'Cañitas De Felipe Pescador'),
    @('B1200', 'El Plateado De Joaquín Amaro'),
    @('B1210', 'Jiménez Del Teul'),
    @('B1216', 'Mezquital Del Oro'),
    @('B1220', 'Moyahua De Estrada'),
    @('B1221', 'Nochistlán De Mejía'),
    @('B1222', 'Noria De Ángeles'),
    @('B1233', 'Teúl De González Ortega'),
    @('B1234', 'Tlaltenango De Sánchez Román'),
    @('B1237', 'Villa De Cos'),
)
foreach ($fix in $stringFixes) {
    $ws.Range($fix[0]).Value = $fix[1]
}

# 3. Floating point precision corrections in column D (percentage values)
$floatFixes = @(
    @('D82', 0.0009323674962346696),
    @('D87', 0.0009323674962346696),
    @('D139', 0.0009323674962346696),
    @('D142', 0.0009323674962346696),
    @('D159', 0.0009323674962346696),
    @('D162', 0.0009323674962346696),
    @('D376', 0.0009323674962346696),
    @('D378', 0.0009323674962346696),
    @('D383', 0.0009323674962346696),
    @('D385', 0.0009323674962346696),
    @('D424', 0.0009323674962346696),
    @('D529', 0.0009323674962346696),
    @('D532', 0.0009323674962346696),
    @('D535', 0.0009323674962346696),
    @('D633', 0.0009323674962346696),
    @('D926', 0.0009323674962346696),
    @('D947', 0.0009323674962346696),
    @('D1021', 0.0009323674962346696),
    @('D1185', 0.0009323674962346696),
    @('D1217', 0.0009323674962346696),
    @('D361', 0.009682277845513875),
)
foreach ($fix in $floatFixes) {
    $ws.Range($fix[0]).Value = $fix[1]
}

# 4. Remove obsolete footer/metadata rows 1246-1250
$ws.Range("A1246:D1250").EntireRow.Delete()
